# Apply edits described by the diff to before.xlsx

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Wafa Al Shamsi"
$wsSummary.Range("B4").Value = 2767.2
$wsSummary.Range("B6").Value = 4654
$wsSummary.Range("B7").Value = 10318
$wsSummary.Range("B8").Value = -5664
$wsSummary.Range("B9").Value = 0.45

# --- Assets sheet ---
$wsAssets = $wb.Worksheets.Item("Assets")
# Update row 2 in place
$wsAssets.Range("A2").Value = "Liquid Assets"
$wsAssets.Range("B2").Value = "Savings Account"
$wsAssets.Range("C2").Value = 4654
# Delete rows 3 and 4 (old "Vehicles/Mid-range Car" and old "Liquid Assets/Savings Account" rows),
# which shifts the TOTAL ASSETS row (formerly row 5) up to row 3.
$wsAssets.Rows("3:4").Delete()
# Update the new TOTAL ASSETS row (now row 3)
$wsAssets.Range("C3").Value = 4654

# --- Liabilities sheet ---
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
# Update row 2 in place
$wsLiabilities.Range("A2").Value = "Credit Cards"
$wsLiabilities.Range("B2").Value = "Credit Card Balance"
$wsLiabilities.Range("C2").Value = 10318
$wsLiabilities.Range("D2").Value = 516
$wsLiabilities.Range("E2").Value = 1
# Delete rows 3 and 4 (old "Personal Loans" and old "Credit Cards" rows),
# which shifts the TOTAL LIABILITIES row (formerly row 5) up to row 3.
$wsLiabilities.Rows("3:4").Delete()
# Update the new TOTAL LIABILITIES row (now row 3)
$wsLiabilities.Range("C3").Value = 10318
